$d = $word.ActiveDocument

$replacements = @(
    @("2023-07-17 Monday", "2023-07-18 Tuesday"),
    @("30-26=", "31+37="),
    @("52+35=", "95-45="),
    @("80-46=", "44-15="),
    @("73-49=", "77-44="),
    @("83+7=", "59-14="),
    @("74-28=", "2+35="),
    @("25-23=", "6+88="),
    @("14+12=", "1+68="),
    @("40+10=", "89-49="),
    @("57-42=", "53+18="),
    @("28+44=", "97-53="),
    @("10+74=", "69-1="),
    @("95-68=", "67+4="),
    @("72-25=", "93-89="),
    @("24+21=", "37+33="),
    @("36+13=", "79-75="),
    @("26-12=", "50-47="),
    @("19+72=", "33+12="),
    @("59-12=", "51-3="),
    @("80-36=", "85-57="),
    @("83-57=", "70-9="),
    @("17+1=", "9+24="),
    @("41+5=", "52-41="),
    @("63-54=", "29+4="),
    @("80-41=", "52-7="),
    @("21+62=", "40+48="),
    @("71+24=", "28+68="),
    @("68-39=", "72+16="),
    @("32+61=", "36+22="),
    @("44+7=", "30+67="),
    @("80-45=", "14+39="),
    @("42+43=", "71-64="),
    @("45+0=", "25-0="),
    @("76-58=", "37+8="),
    @("3+83=", "2+3="),
    @("31+45=", "40+23="),
    @("93-3=", "77+14="),
    @("72-5=", "34-22="),
    @("86-15=", "75-68="),
    @("51+16=", "96-77="),
    @("54+0=", "35-17="),
    @("55-48=", "5+77="),
    @("85-60=", "63-39="),
    @("65+5=", "20+35="),
    @("11+51=", "22+6="),
    @("86-1=", "98-80="),
    @("66-63=", "82+1="),
    @("78-27=", "45+1="),
    @("81-46=", "17+68="),
    @("97-87=", "23+47="),
    @("52+32=", "29+64="),
    @("1+75=", "96-57="),
    @("20+16=", "8+47="),
    @("28+36=", "39-2="),
    @("95-42=", "73-57="),
    @("39+11=", "36+9="),
    @("71-53=", "7+22="),
    @("6+49=", "42+55="),
    @("92-83=", "18+13="),
    @("33-4=", "77+3="),
    @("68-5=", "23-22="),
    @("97-71=", "99-39="),
    @("82-21=", "82-28="),
    @("66+5=", "9+4="),
    @("86-25=", "46-2="),
    @("30+47=", "36+1="),
    @("1+89=", "97-95="),
    @("35-0=", "78+17="),
    @("80-12=", "70-19="),
    @("81-48=", "74+20="),
    @("64+16=", "56-47="),
    @("19+58=", "89-62="),
    @("84-53=", "49-42="),
    @("45-31=", "26-16="),
    @("86-7=", "25+24="),
    @("27+67=", "59+36="),
    @("33+4=", "69+17="),
    @("93-79=", "15+71="),
    @("35-16=", "0+89="),
    @("77+6=", "79-54="),
    @("71-54=", "85-51="),
    @("90-13=", "2+81="),
    @("7+70=", "68+4="),
    @("32+38=", "99-93="),
    @("80+4=", "45+9="),
    @("37+17=", "67+32="),
    @("60-53=", "46+26="),
    @("3+10=", "57+33="),
    @("34+52=", "15+28="),
    @("83-26=", "71-7="),
    @("17-4=", "69-51="),
    @("5+91=", "90-19="),
    @("99-23=", "6+86="),
    @("81-45=", "53-7="),
    @("13+42=", "93-23="),
    @("74-11=", "62+6="),
    @("40+11=", "90-70="),
    @("51+5=", "42+41="),
    @("29+7=", "6+19="),
    @("51+17=", "54-6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
